$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# Row 16 -> D22
$ws.Range("B16").Value = "D22"
$ws.Range("C15").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = 43856
$ws.Range("D16").Value = "Started working on Intel Image classification problem. Also watched some videos about CNN"

# Row 17 -> D23
$ws.Range("B17").Value = "D23"
$ws.Range("C15").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 43857
$ws.Range("D17").Value = "Out sick"

# Row 18 -> D24
$ws.Range("B18").Value = "D24"
$ws.Range("C15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 43858
$ws.Range("D18").Value = "Completed lessons 6.11 to 6.24"
